# Update gh-pages to output generated at 456a3b4
# Applies the same set of data refreshes to both the "展览" and "全部类型"
# worksheets (they carry duplicate copies of the exhibition listing).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 4: "合肥·二次元晚会" event — rename, bump interest count, refresh cover image
    $ws.Range("C4").Value = "合肥·二次元晚会（免费活动）"
    $ws.Range("F4").Value = 54
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202404/dRMsa2dn1713437265983.jpeg"

    # Remaining rows only have their "想去人数" (interest count) refreshed
    $ws.Range("F5").Value = 363
    $ws.Range("F6").Value = 11029
    $ws.Range("F7").Value = 518
    $ws.Range("F8").Value = 96
    $ws.Range("F11").Value = 145
    $ws.Range("F18").Value = 308
    $ws.Range("F19").Value = 1165
    $ws.Range("F21").Value = 883
    $ws.Range("F22").Value = 107
}
